$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "no.ID"
$ws.Range("B16").Value = "no.ID"
$ws.Range("C16").Value = "no.ID"

$ws.Range("A17").Value = "Noise"
$ws.Range("B17").Value = "Noise"
$ws.Range("C17").Value = "Noise"

$ws.Range("C18").Select()
